$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 3.156457666666666
$ws.Cells.Item(2, 8).Value = 9.469372999999999
$ws.Cells.Item(2, 9).Value = 0.03431910465203977
$ws.Cells.Item(2, 10).Value = 0.03431910465203977
$ws.Cells.Item(2, 13).Value = 4.277890333333334
$ws.Cells.Item(2, 14).Value = 12.833671
$ws.Cells.Item(2, 15).Value = 0.04123357425337639
$ws.Cells.Item(2, 16).Value = 0.04123357425337638
$ws.Cells.Item(2, 17).Value = 13.50297973980923
$ws.Cells.Item(2, 18).Value = 121.526817658283
$ws.Cells.Item(2, 19).Value = 0.001415099349979277
$ws.Cells.Item(2, 20).Value = 0.001415099349979277
$ws.Cells.Item(3, 7).Value = 3.156457666666666
$ws.Cells.Item(3, 8).Value = 9.469372999999999
$ws.Cells.Item(3, 9).Value = 0.03431910465203977
$ws.Cells.Item(3, 10).Value = 0.03431910465203977
$ws.Cells.Item(3, 15).Value = 0.4451428460610328
$ws.Cells.Item(3, 16).Value = 0.4451428460610327
$ws.Cells.Item(3, 17).Value = 145.7733155691919
$ws.Cells.Item(3, 18).Value = 1311.959840122727
$ws.Cells.Item(3, 19).Value = 0.01527690391907541
$ws.Cells.Item(3, 20).Value = 0.01527690391907541
$ws.Cells.Item(4, 7).Value = 3.156457666666666
$ws.Cells.Item(4, 8).Value = 9.469372999999999
$ws.Cells.Item(4, 9).Value = 0.03431910465203977
$ws.Cells.Item(4, 10).Value = 0.03431910465203977
$ws.Cells.Item(4, 13).Value = 8.558147333333332
$ws.Cells.Item(4, 14).Value = 25.674442
$ws.Cells.Item(4, 15).Value = 0.08248996024761777
$ws.Cells.Item(4, 16).Value = 0.08248996024761777
$ws.Cells.Item(4, 17).Value = 27.01342976276288
$ws.Cells.Item(4, 18).Value = 243.120867864866
$ws.Cells.Item(4, 19).Value = 0.002830981578480595
$ws.Cells.Item(4, 20).Value = 0.002830981578480595
$ws.Cells.Item(5, 7).Value = 3.156457666666666
$ws.Cells.Item(5, 8).Value = 9.469372999999999
$ws.Cells.Item(5, 9).Value = 0.03431910465203977
$ws.Cells.Item(5, 10).Value = 0.03431910465203977
$ws.Cells.Item(5, 13).Value = 44.72914066666667
$ws.Cells.Item(5, 14).Value = 134.187422
$ws.Cells.Item(5, 15).Value = 0.4311336194379731
$ws.Cells.Item(5, 16).Value = 0.431133619437973
$ws.Cells.Item(5, 17).Value = 141.1856389807118
$ws.Cells.Item(5, 18).Value = 1270.670750826406
$ws.Cells.Item(5, 19).Value = 0.01479611980450449
$ws.Cells.Item(5, 20).Value = 0.01479611980450448
$ws.Cells.Item(6, 9).Value = 0.6977360249749448
$ws.Cells.Item(6, 10).Value = 0.6977360249749449
$ws.Cells.Item(6, 13).Value = 4.277890333333334
$ws.Cells.Item(6, 14).Value = 12.833671
$ws.Cells.Item(6, 15).Value = 0.04123357425337639
$ws.Cells.Item(6, 16).Value = 0.04123357425337638
$ws.Cells.Item(6, 17).Value = 274.5268416672325
$ws.Cells.Item(6, 18).Value = 2470.741575005092
$ws.Cells.Item(6, 19).Value = 0.02877015019506007
$ws.Cells.Item(6, 20).Value = 0.02877015019506007
$ws.Cells.Item(7, 9).Value = 0.6977360249749448
$ws.Cells.Item(7, 10).Value = 0.6977360249749449
$ws.Cells.Item(7, 15).Value = 0.4451428460610328
$ws.Cells.Item(7, 16).Value = 0.4451428460610327
$ws.Cells.Item(7, 19).Value = 0.3105921999566588
$ws.Cells.Item(7, 20).Value = 0.3105921999566588
$ws.Cells.Item(8, 9).Value = 0.6977360249749448
$ws.Cells.Item(8, 10).Value = 0.6977360249749449
$ws.Cells.Item(8, 13).Value = 8.558147333333332
$ws.Cells.Item(8, 14).Value = 25.674442
$ws.Cells.Item(8, 15).Value = 0.08248996024761777
$ws.Cells.Item(8, 16).Value = 0.08248996024761777
$ws.Cells.Item(8, 17).Value = 549.2055604221536
$ws.Cells.Item(8, 18).Value = 4942.850043799383
$ws.Cells.Item(8, 19).Value = 0.05755621696351404
$ws.Cells.Item(8, 20).Value = 0.05755621696351405
$ws.Cells.Item(9, 9).Value = 0.6977360249749448
$ws.Cells.Item(9, 10).Value = 0.6977360249749449
$ws.Cells.Item(9, 13).Value = 44.72914066666667
$ws.Cells.Item(9, 14).Value = 134.187422
$ws.Cells.Item(9, 15).Value = 0.4311336194379731
$ws.Cells.Item(9, 16).Value = 0.431133619437973
$ws.Cells.Item(9, 17).Value = 2870.421810963371
$ws.Cells.Item(9, 18).Value = 25833.79629867034
$ws.Cells.Item(9, 19).Value = 0.3008174578597119
$ws.Cells.Item(9, 20).Value = 0.3008174578597119
$ws.Cells.Item(10, 7).Value = 21.527469
$ws.Cells.Item(10, 8).Value = 64.582407
$ws.Cells.Item(10, 9).Value = 0.2340609441104101
$ws.Cells.Item(10, 10).Value = 0.2340609441104101
$ws.Cells.Item(10, 13).Value = 4.277890333333334
$ws.Cells.Item(10, 14).Value = 12.833671
$ws.Cells.Item(10, 15).Value = 0.04123357425337639
$ws.Cells.Item(10, 16).Value = 0.04123357425337638
$ws.Cells.Item(10, 17).Value = 92.09215153623302
$ws.Cells.Item(10, 18).Value = 828.8293638260972
$ws.Cells.Item(10, 19).Value = 0.009651169318791975
$ws.Cells.Item(10, 20).Value = 0.009651169318791973
$ws.Cells.Item(11, 7).Value = 21.527469
$ws.Cells.Item(11, 8).Value = 64.582407
$ws.Cells.Item(11, 9).Value = 0.2340609441104101
$ws.Cells.Item(11, 10).Value = 0.2340609441104101
$ws.Cells.Item(11, 15).Value = 0.4451428460610328
$ws.Cells.Item(11, 16).Value = 0.4451428460610327
$ws.Cells.Item(11, 17).Value = 994.193765081277
$ws.Cells.Item(11, 18).Value = 8947.743885731494
$ws.Cells.Item(11, 19).Value = 0.1041905548130403
$ws.Cells.Item(11, 20).Value = 0.1041905548130403
$ws.Cells.Item(12, 7).Value = 21.527469
$ws.Cells.Item(12, 8).Value = 64.582407
$ws.Cells.Item(12, 9).Value = 0.2340609441104101
$ws.Cells.Item(12, 10).Value = 0.2340609441104101
$ws.Cells.Item(12, 13).Value = 8.558147333333332
$ws.Cells.Item(12, 14).Value = 25.674442
$ws.Cells.Item(12, 15).Value = 0.08248996024761777
$ws.Cells.Item(12, 16).Value = 0.08248996024761777
$ws.Cells.Item(12, 17).Value = 184.235251415766
$ws.Cells.Item(12, 18).Value = 1658.117262741894
$ws.Cells.Item(12, 19).Value = 0.01930767797518761
$ws.Cells.Item(12, 20).Value = 0.01930767797518761
$ws.Cells.Item(13, 7).Value = 21.527469
$ws.Cells.Item(13, 8).Value = 64.582407
$ws.Cells.Item(13, 9).Value = 0.2340609441104101
$ws.Cells.Item(13, 10).Value = 0.2340609441104101
$ws.Cells.Item(13, 13).Value = 44.72914066666667
$ws.Cells.Item(13, 14).Value = 134.187422
$ws.Cells.Item(13, 15).Value = 0.4311336194379731
$ws.Cells.Item(13, 16).Value = 0.431133619437973
$ws.Cells.Item(13, 17).Value = 962.905189098306
$ws.Cells.Item(13, 18).Value = 8666.146701884754
$ws.Cells.Item(13, 19).Value = 0.1009115420033902
$ws.Cells.Item(13, 20).Value = 0.1009115420033902
$ws.Cells.Item(14, 7).Value = 3.116432666666667
$ws.Cells.Item(14, 8).Value = 9.349298000000001
$ws.Cells.Item(14, 9).Value = 0.03388392626260537
$ws.Cells.Item(14, 10).Value = 0.03388392626260537
$ws.Cells.Item(14, 13).Value = 4.277890333333334
$ws.Cells.Item(14, 14).Value = 12.833671
$ws.Cells.Item(14, 15).Value = 0.04123357425337639
$ws.Cells.Item(14, 16).Value = 0.04123357425337638
$ws.Cells.Item(14, 17).Value = 13.33175717921756
$ws.Cells.Item(14, 18).Value = 119.985814612958
$ws.Cells.Item(14, 19).Value = 0.001397155389545069
$ws.Cells.Item(14, 20).Value = 0.001397155389545069
$ws.Cells.Item(15, 7).Value = 3.116432666666667
$ws.Cells.Item(15, 8).Value = 9.349298000000001
$ws.Cells.Item(15, 9).Value = 0.03388392626260537
$ws.Cells.Item(15, 10).Value = 0.03388392626260537
$ws.Cells.Item(15, 15).Value = 0.4451428460610328
$ws.Cells.Item(15, 16).Value = 0.4451428460610327
$ws.Cells.Item(15, 17).Value = 143.9248583517002
$ws.Cells.Item(15, 18).Value = 1295.323725165302
$ws.Cells.Item(15, 19).Value = 0.01508318737225833
$ws.Cells.Item(15, 20).Value = 0.01508318737225833
$ws.Cells.Item(16, 7).Value = 3.116432666666667
$ws.Cells.Item(16, 8).Value = 9.349298000000001
$ws.Cells.Item(16, 9).Value = 0.03388392626260537
$ws.Cells.Item(16, 10).Value = 0.03388392626260537
$ws.Cells.Item(16, 13).Value = 8.558147333333332
$ws.Cells.Item(16, 14).Value = 25.674442
$ws.Cells.Item(16, 15).Value = 0.08248996024761777
$ws.Cells.Item(16, 16).Value = 0.08248996024761777
$ws.Cells.Item(16, 17).Value = 26.67088991574622
$ws.Cells.Item(16, 18).Value = 240.038009241716
$ws.Cells.Item(16, 19).Value = 0.002795083730435529
$ws.Cells.Item(16, 20).Value = 0.002795083730435529
$ws.Cells.Item(17, 7).Value = 3.116432666666667
$ws.Cells.Item(17, 8).Value = 9.349298000000001
$ws.Cells.Item(17, 9).Value = 0.03388392626260537
$ws.Cells.Item(17, 10).Value = 0.03388392626260537
$ws.Cells.Item(17, 13).Value = 44.72914066666667
$ws.Cells.Item(17, 14).Value = 134.187422
$ws.Cells.Item(17, 15).Value = 0.4311336194379731
$ws.Cells.Item(17, 16).Value = 0.431133619437973
$ws.Cells.Item(17, 17).Value = 139.3953551255285
$ws.Cells.Item(17, 18).Value = 1254.558196129756
$ws.Cells.Item(17, 19).Value = 0.01460849977036644
$ws.Cells.Item(17, 20).Value = 0.01460849977036644